# Made changes to run the test case for gmail login.
$wb = $excel.ActiveWorkbook

$wsRun = $wb.Worksheets.Item("RunManager")
$wsData = $wb.Worksheets.Item("TestData")

# --- Sheet "TestData" (sheet2): refresh the reusable data-driven values ---
$wsData.Range("D2").Value = "checkWhetherGlobalViewIsSelectedByDefault"
$wsData.Range("D3").Value = "validateEnteringShipmentID"
$wsData.Range("D4").Value = "test11"
$wsData.Range("D5").Value = "test12"
$wsData.Range("D6").Value = "test2"
$wsData.Range("D7").Value = "test3"

# --- Sheet "RunManager" (sheet1) ---
# Only run the "loginTest" test case: remove the other two rows (test2, test3)
# and rename the remaining test case from "test1" to "loginTest".
$wsRun.Rows.Item(4).Delete() | Out-Null
$wsRun.Rows.Item(3).Delete() | Out-Null
$wsRun.Range("A2").Value = "loginTest"

# --- Sheet "TestData" (sheet2) ---
# Update the test data row for the login test case and reuse it for the
# gmail login scenario, wiring the username cell up as a mailto hyperlink.
$wsData.Range("A2").Value = "loginTest"
$wsData.Range("B2").Value = "g.amaresh18@gmail.com"
$wsData.Hyperlinks.Add($wsData.Range("B2"), "mailto:g.amaresh18@gmail.com") | Out-Null
$wsData.Range("F3").Select() | Out-Null

# Restore "RunManager" as the active sheet/tab with its own selection, since
# that is the sheet that was active/selected in the saved workbook.
$wsRun.Range("B2").Select() | Out-Null
